$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback DateTime for the
# 52032f82-... row (row 3)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-12 18:35:00"
$wsZh.Range("H3").Value = "2016-03-12 18:35:21"

# de-de sheet: same row/columns
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-12 18:35:04"
$wsDe.Range("H3").Value = "2016-03-12 18:35:26"
